$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("4-6 Pour Over Recipe")

# Update C17: was formula =2/3, now plain value 0.5
$ws.Range("C17").Value = 0.5

# Update C18: was formula =1/3, now plain value 0.5
$ws.Range("C18").Value = 0.5

# Update C19: was plain value 0, now blank
$ws.Range("C19").ClearContents()

# Break the shared-formula group E16:E19 into individual formulas
$ws.Range("E16").Formula = "=E15+D16"
$ws.Range("E17").Formula = "=E16+D17"
$ws.Range("E18").Formula = "=E17+D18"
$ws.Range("E19").Formula = "=E18+D19"

# Remove the data validation rule on C17:C19 (the "60% of Water" rule)
$ws.Range("C17:C19").Validation.Delete()

# Update the active selection to C17
$ws.Range("C17").Select()
